$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.08511433333333333
$ws.Range("H2").Value = 0.255343
$ws.Range("I2").Value = 0.3100414411767206
$ws.Range("J2").Value = 0.3100414411767206
$ws.Range("M2").Value = 0.003615
$ws.Range("N2").Value = 0.010845
$ws.Range("Q2").Value = 0.000307688315
$ws.Range("R2").Value = 0.002769194835
$ws.Range("S2").Value = 0.3100414411767206
$ws.Range("T2").Value = 0.3100414411767206

$ws.Range("I3").Value = 0.01968121984951013
$ws.Range("J3").Value = 0.01968121984951013
$ws.Range("M3").Value = 0.003615
$ws.Range("N3").Value = 0.010845
$ws.Range("Q3").Value = 0.000019531845
$ws.Range("R3").Value = 0.000175786605
$ws.Range("S3").Value = 0.01968121984951013
$ws.Range("T3").Value = 0.01968121984951013

$ws.Range("G4").Value = 0.1840083333333333
$ws.Range("H4").Value = 0.552025
$ws.Range("I4").Value = 0.6702773389737693
$ws.Range("J4").Value = 0.6702773389737693
$ws.Range("M4").Value = 0.003615
$ws.Range("N4").Value = 0.010845
$ws.Range("Q4").Value = 0.000665190125
$ws.Range("R4").Value = 0.005986711125
$ws.Range("S4").Value = 0.6702773389737693
$ws.Range("T4").Value = 0.6702773389737693
